# Applies the cryptocurrency market-data refresh described in the commit:
# "Updated cryptos list on Fri Dec  8 23:20:38 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when it looks numeric
# (e.g. "238.70"), by forcing a Text number format before the write and
# clearing back to the default cell style afterwards so no formatting
# change is left behind - only the displayed text changes.
function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '44.357.28'
$ws.Range('E2').Value = '  +2.45%  '
$ws.Range('D3').Value = '2.364.52'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  +3.99%  '
Set-TextValue 'D6' '238.70'
$ws.Range('E6').Value = '  +2.74%  '
Set-TextValue 'D7' '73.98'
$ws.Range('E7').Value = '  +9.14%  '
$ws.Range('E8').Value = '  -0.02%  '
Set-TextValue 'D9' '0.553'
$ws.Range('E9').Value = '  +21.16%  '
$ws.Range('E10').Value = '  +6.45%  '
Set-TextValue 'D11' '30.34'
$ws.Range('E11').Value = '  +15.39%  '
$ws.Range('E12').Value = '  +2.16%  '
$ws.Range('D13').Value = '2.711.93'
$ws.Range('E13').Value = '  +0.26%  '
Set-TextValue 'D14' '16.93'
$ws.Range('E14').Value = '  +8.26%  '
Set-TextValue 'D15' '6.86'
$ws.Range('E15').Value = '  +10.00%  '
Set-TextValue 'D16' '0.907'
$ws.Range('E16').Value = '  +8.15%  '
$ws.Range('D17').Value = '2.361.24'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').Value = '44.492.85'
$ws.Range('E18').Value = '  +2.84%  '
$ws.Range('E19').Value = '  +4.64%  '
Set-TextValue 'D20' '77.42'
$ws.Range('E20').Value = '  +4.73%  '
Set-TextValue 'D21' '6.49'
$ws.Range('E21').Value = '  +4.16%  '
Set-TextValue 'D22' '255.02'
$ws.Range('E22').Value = '  +2.54%  '
Set-TextValue 'D23' '3.86'
$ws.Range('E23').Value = '  -3.31%  '
Set-TextValue 'D24' '0.999'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  +2.60%  '
Set-TextValue 'D26' '10.37'
$ws.Range('E26').Value = '  +4.80%  '
Set-TextValue 'D27' '2.25'
$ws.Range('E27').Value = '  -1.10%  '
Set-TextValue 'D28' '22.63'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D29' '173.88'
$ws.Range('E29').Value = '  +0.62%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D30' '1.59'
$ws.Range('E30').Value = '  +4.13%  '
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('E32').Value = '  +5.20%  '
Set-TextValue 'D33' '0.0744'
$ws.Range('E33').Value = '  +7.23%  '
Set-TextValue 'D34' '5.20'
$ws.Range('E34').Value = '  +4.00%  '
Set-TextValue 'D35' '5.25'
$ws.Range('E35').Value = '  +3.95%  '
$ws.Range('E36').Value = '  +7.61%  '
Set-TextValue 'D37' '2.44'
$ws.Range('E37').Value = '  -2.32%  '
Set-TextValue 'D38' '6.48'
$ws.Range('E38').Value = '  -0.27%  '
Set-TextValue 'D39' '0.0271'
$ws.Range('E39').Value = '  +6.66%  '
Set-TextValue 'D40' '19.39'
$ws.Range('E40').Value = '  +7.13%  '
$ws.Range('E41').Value = '  -0.06%  '
Set-TextValue 'D42' '8.84'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('E43').Value = '  +3.16%  '
Set-TextValue 'D44' '0.0990'
$ws.Range('E44').Value = '  +4.29%  '
$ws.Range('E45').Value = '  +1.58%  '
$ws.Range('E46').Value = '  +12.59%  '
Set-TextValue 'D47' '99.19'
$ws.Range('E47').Value = '  +0.72%  '
Set-TextValue 'D48' '4.46'
$ws.Range('E48').Value = '  -0.12%  '
Set-TextValue 'D49' '2.39'
$ws.Range('E49').Value = '  +5.38%  '
$ws.Range('D50').Value = '1.444.88'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').Value = '2.585.39'
$ws.Range('E51').Value = '  +0.23%  '
